$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BPTBfRN")

# Rename the existing "hydrogen" plant type to "hydrogen combustion turbine"
$ws.Cells.Item(24, 1).Value2 = "hydrogen combustion turbine"

# Add a new row for "hydrogen combined cycle" as a second hydrogen-based plant type
$ws.Cells.Item(25, 1).Value2 = "hydrogen combined cycle"
$ws.Cells.Item(25, 2).Value2 = 1

# Give the new/renamed hydrogen rows a distinct font (black, vertically centered)
$rng24 = $ws.Range("A24")
$rng24.Font.Color = 0
$rng24.VerticalAlignment = -4108

# Copy the resolved formatting onto the new row instead of re-deriving it,
# so we don't leave a stray intermediate style behind
$rng24.Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection in BPTBfRN down past the new data, then bring the
# "About" sheet to the front as the active tab
[void]$ws.Range("B30").Select()
[void]$wsAbout.Select()

"done"
